$d = $word.ActiveDocument

# 1) Add a trailing period to the first paragraph's text.
$d.Content.Find.Execute("Este es un texto", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Este es un texto.", 2)

# 2) Add a trailing period to the second paragraph's second run text.
$d.Content.Find.Execute("Es la continuación de la misma segunda linea con otros estilos para seguir probando", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Es la continuación de la misma segunda linea con otros estilos para seguir probando.", 2)

# 3) Add a trailing period to the third paragraph's text.
$d.Content.Find.Execute("Este es el estilo correspondiente al tercer párrafo del documento que se esta creando", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Este es el estilo correspondiente al tercer párrafo del documento que se esta creando.", 2)

# 4) Create a new fourth paragraph with its own styles (right-aligned,
#    Verdana, blue, 13pt/26 half-points).
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newRange = $newPara.Range
$newRange.Text = "Este es el siguiente párrafo, (párrafo 4) creado con los diferentes estilos ya predefinidos"

# Restrict formatting to the text itself (exclude the paragraph mark) so
# no stray rPr ends up on the pPr.
$textRange = $d.Range($newRange.Start, $newRange.End - 1)
$textRange.Font.Name = "Verdana"
$textRange.Font.Color = 16711680
$textRange.Font.Size = 13

$newPara.Alignment = 2
